$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "EkDhoTeenChar"
$ws.Range("C2").Select()
